$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value for the "dSF" (column F) updates described in the diff.
$changes = @{
    "F2"  = -10
    "F3"  = -6
    "F4"  = -6
    "F6"  = -4
    "F9"  = -4
    "F11" = -4
    "F12" = 6
    "F19" = -3
    "F21" = 3
    "F30" = -4
    "F33" = 5
    "F37" = -4
    "F38" = 2
    "F40" = -3
    "F45" = 6
    "F48" = -4
    "F49" = 10
    "F53" = 3
    "F56" = 1
}

foreach ($cell in $changes.Keys) {
    $ws.Range($cell).Value = $changes[$cell]
}
